$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matching source data)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "61.510.89"
$ws.Range("E2").Value = "  -1.40%  "

$ws.Range("D3").Value = "2.442.95"
$ws.Range("E3").Value = "  -0.44%  "

$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "576.14"
$ws.Range("E5").Value = "  -0.87%  "

$ws.Range("D6").Value = "140.62"
$ws.Range("E6").Value = "  -2.18%  "

$ws.Range("E7").Value = "  +0.15%  "

$ws.Range("D8").Value = "0.531"
$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("D9").Value = "2.433.36"
$ws.Range("E9").Value = "  -0.77%  "

$ws.Range("D10").Value = "0.109"
$ws.Range("E10").Value = "  +2.01%  "

$ws.Range("E11").Value = "  +1.55%  "

$ws.Range("D12").Value = "5.16"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("D13").Value = "0.340"
$ws.Range("E13").Value = "  -1.63%  "

$ws.Range("D14").Value = "26.06"
$ws.Range("E14").Value = "  -1.62%  "

$ws.Range("D15").Value = "2.894.25"
$ws.Range("E15").Value = "  +1.02%  "

$ws.Range("D16").Value = "0.0000170"
$ws.Range("E16").Value = "  -1.27%  "

$ws.Range("D17").Value = "61.538.71"
$ws.Range("E17").Value = "  -1.20%  "

$ws.Range("D18").Value = "2.435.20"
$ws.Range("E18").Value = "  -0.27%  "

$ws.Range("E19").Value = "  -3.46%  "

$ws.Range("D20").Value = "7.27"
$ws.Range("E20").Value = "  +1.85%  "

$ws.Range("D21").Value = "324.95"
$ws.Range("E21").Value = "  -1.73%  "

$ws.Range("E22").Value = "  -1.76%  "

$ws.Range("D23").Value = "6.01"
$ws.Range("E23").Value = "  +0.19%  "

$ws.Range("E24").Value = "  -0.06%  "

$ws.Range("D25").Value = "1.91"
$ws.Range("E25").Value = "  -2.34%  "

$ws.Range("D26").Value = "64.87"
$ws.Range("E26").Value = "  -1.26%  "

$ws.Range("D27").Value = "9.11"
$ws.Range("E27").Value = "  -3.20%  "

$ws.Range("D28").Value = "579.78"
$ws.Range("E28").Value = "  -7.93%  "

$ws.Range("D29").Value = "2.569.80"
$ws.Range("E29").Value = "  +0.37%  "

$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("D31").Value = "0.0₃0917"
$ws.Range("E31").Value = "  -4.02%  "

$ws.Range("D32").Value = "7.89"
$ws.Range("E32").Value = "  -1.87%  "

$ws.Range("D33").Value = "1.35"
$ws.Range("E33").Value = "  -5.31%  "

$ws.Range("E34").Value = "  -1.20%  "

$ws.Range("E35").Value = "  -6.29%  "

$ws.Range("E36").Value = "  +0.28%  "

$ws.Range("D37").Value = "4.66"
$ws.Range("E37").Value = "  -5.87%  "

$ws.Range("D38").Value = "0.371"
$ws.Range("E38").Value = "  -1.24%  "

$ws.Range("D39").Value = "151.84"
$ws.Range("E39").Value = "  +1.11%  "

$ws.Range("E40").Value = "  -3.77%  "

$ws.Range("D41").Value = "18.33"
$ws.Range("E41").Value = "  -0.20%  "

$ws.Range("D42").Value = "5.12"
$ws.Range("E42").Value = "  -2.90%  "

$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "41.70"
$ws.Range("E44").Value = "  -2.63%  "

$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "1.67"
$ws.Range("E45").Value = "  -5.51%  "

$ws.Range("E46").Value = "  -5.90%  "

$ws.Range("D47").Value = "0.0₆0293"
$ws.Range("E47").Value = "  +24.46%  "

$ws.Range("D48").Value = "142.44"
$ws.Range("E48").Value = "  -0.71%  "

$ws.Range("D49").Value = "3.55"
$ws.Range("E49").Value = "  -2.61%  "

$ws.Range("D50").Value = "0.596"
$ws.Range("E50").Value = "  -0.90%  "

$ws.Range("D51").Value = "0.0509"
$ws.Range("E51").Value = "  -3.22%  "

# Restore original (default) cell formatting/style for column D
$ws.Range("D2:D51").ClearFormats()
